# Scheduled runner update: refresh market-board price snapshots
# (currentAveragePrice / LevePrice / LeveProfit columns) across all
# gathering-class profit sheets. Generated from the latest pull.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: ALC
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")
# Leve row 86
$ws.Range("H86").Value = 5038.923
$ws.Range("I86").Value = 2000
$ws.Range("J86").Value = 5591.4546
$ws.Range("K86").Value = 2000
$ws.Range("L86").Value = 5591.4546
$ws.Range("M86").Value = -877
$ws.Range("N86").Value = -7837.4546
# Leve row 89
$ws.Range("H89").Value = 5038.923
$ws.Range("I89").Value = 2000
$ws.Range("J89").Value = 5591.4546
$ws.Range("K89").Value = 10000
$ws.Range("L89").Value = 27957.273
$ws.Range("M89").Value = -4384
$ws.Range("N89").Value = -39189.273
# Leve row 100
$ws.Range("H100").Value = 7237.357
$ws.Range("I100").Value = 3899.8333
$ws.Range("J100").Value = 9740.5
$ws.Range("K100").Value = 3899.8333
$ws.Range("L100").Value = 9740.5
$ws.Range("M100").Value = -3358.8333
$ws.Range("N100").Value = -10822.5
# Leve row 107
$ws.Range("H107").Value = 3000
$ws.Range("I107").Value = 3000
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3000
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1080
$ws.Range("N107").ClearContents()
# Leve row 112
$ws.Range("H112").Value = 1419.6451
$ws.Range("J112").Value = 1481.2142
$ws.Range("L112").Value = 4443.642599999999
$ws.Range("N112").Value = -6659.642599999999
# Leve row 116
$ws.Range("H116").Value = 38091
$ws.Range("I116").Value = 54012.383
$ws.Range("K116").Value = 54012.383
$ws.Range("M116").Value = -50570.383
# Leve row 127
$ws.Range("H127").Value = 5662.25
$ws.Range("I127").Value = 5662.25
$ws.Range("K127").Value = 16986.75
$ws.Range("M127").Value = -12026.75
# Leve row 132
$ws.Range("H132").Value = 3256.7222
$ws.Range("I132").Value = 3213.0588
$ws.Range("K132").Value = 9639.1764
$ws.Range("M132").Value = -7109.1764
# Leve row 135
$ws.Range("H135").Value = 1297.8286
$ws.Range("I135").Value = 1295.9678
$ws.Range("K135").Value = 11663.7102
$ws.Range("M135").Value = -9128.7102
# Leve row 141
$ws.Range("H141").Value = 1074
$ws.Range("I141").Value = 283.5
$ws.Range("J141").Value = 2655
$ws.Range("K141").Value = 850.5
$ws.Range("L141").Value = 7965
$ws.Range("M141").Value = 4329.5
$ws.Range("N141").Value = -18325

# ---------------------------------------------------------------
# Sheet: ARM
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")
# Leve row 74
$ws.Range("H74").Value = 5020.3516
$ws.Range("I74").Value = 4622.9062
$ws.Range("K74").Value = 4622.9062
$ws.Range("M74").Value = -3748.9062
# Leve row 77
$ws.Range("H77").Value = 5020.3516
$ws.Range("I77").Value = 4622.9062
$ws.Range("K77").Value = 23114.531
$ws.Range("M77").Value = -18746.531
# Leve row 88
$ws.Range("H88").Value = 1270.6666
$ws.Range("I88").Value = 1020.3571
$ws.Range("J88").Value = 1771.2858
$ws.Range("K88").Value = 1020.3571
$ws.Range("L88").Value = 1771.2858
$ws.Range("M88").Value = -614.3570999999999
$ws.Range("N88").Value = -2583.2858
# Leve row 91
$ws.Range("H91").Value = 1270.6666
$ws.Range("I91").Value = 1020.3571
$ws.Range("J91").Value = 1771.2858
$ws.Range("K91").Value = 1020.3571
$ws.Range("L91").Value = 1771.2858
$ws.Range("M91").Value = 383.6429000000001
$ws.Range("N91").Value = -4579.2858
# Leve row 110
$ws.Range("H110").Value = 741.6842
$ws.Range("I110").Value = 776.41174
$ws.Range("J110").Value = 446.5
$ws.Range("K110").Value = 776.41174
$ws.Range("L110").Value = 446.5
$ws.Range("M110").Value = 1268.58826
$ws.Range("N110").Value = -4536.5

# ---------------------------------------------------------------
# Sheet: BSM
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")
# Leve row 20
$ws.Range("H20").Value = 2014.4706
$ws.Range("I20").Value = 2312
$ws.Range("K20").Value = 2312
$ws.Range("M20").Value = -2065
# Leve row 94
$ws.Range("H94").Value = 899.4
$ws.Range("I94").Value = 948.619
$ws.Range("J94").Value = 784.55554
$ws.Range("K94").Value = 948.619
$ws.Range("L94").Value = 784.55554
$ws.Range("M94").Value = -497.619
$ws.Range("N94").Value = -1686.55554
# Leve row 140
$ws.Range("H140").Value = 183816.75
$ws.Range("J140").Value = 183816.75
$ws.Range("L140").Value = 183816.75
$ws.Range("N140").Value = -194176.75

# ---------------------------------------------------------------
# Sheet: CRP
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")
# Leve row 16
$ws.Range("H16").Value = 752.9167
$ws.Range("I16").Value = 639.5454999999999
$ws.Range("K16").Value = 639.5454999999999
$ws.Range("M16").Value = -352.5454999999999
# Leve row 64
$ws.Range("H64").Value = 68747.25
$ws.Range("J64").Value = 74996.336
$ws.Range("L64").Value = 74996.336
$ws.Range("N64").Value = -75492.336
# Leve row 67
$ws.Range("H67").Value = 68747.25
$ws.Range("J67").Value = 74996.336
$ws.Range("L67").Value = 74996.336
$ws.Range("N67").Value = -76712.336
# Leve row 113
$ws.Range("H113").Value = 752.9167
$ws.Range("I113").Value = 639.5454999999999
$ws.Range("K113").Value = 639.5454999999999
$ws.Range("M113").Value = 1530.4545

# ---------------------------------------------------------------
# Sheet: CUL
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")
# Leve row 87
$ws.Range("H87").Value = 17543.54
$ws.Range("I87").Value = 10500.333
$ws.Range("K87").Value = 31500.999
$ws.Range("M87").Value = -30252.999
# Leve row 88
$ws.Range("H88").Value = 69955
$ws.Range("J88").Value = 69955
$ws.Range("L88").Value = 209865
$ws.Range("N88").Value = -210721
# Leve row 90
$ws.Range("H90").Value = 17543.54
$ws.Range("I90").Value = 10500.333
$ws.Range("K90").Value = 94502.997
$ws.Range("M90").Value = -88262.997
# Leve row 91
$ws.Range("H91").Value = 69955
$ws.Range("J91").Value = 69955
$ws.Range("L91").Value = 209865
$ws.Range("N91").Value = -212829
# Leve row 132
$ws.Range("H132").Value = 2665.8823
$ws.Range("I132").Value = 2667.5
$ws.Range("K132").Value = 24007.5
$ws.Range("M132").Value = -21477.5

# ---------------------------------------------------------------
# Sheet: GSM
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")
# Leve row 82
$ws.Range("H82").Value = 124473.14
$ws.Range("J82").Value = 120218.664
$ws.Range("L82").Value = 120218.664
$ws.Range("N82").Value = -120984.664
# Leve row 85
$ws.Range("H85").Value = 124473.14
$ws.Range("J85").Value = 120218.664
$ws.Range("L85").Value = 120218.664
$ws.Range("N85").Value = -122870.664
# Leve row 113
$ws.Range("H113").Value = 3570.077
$ws.Range("I113").Value = 1679.8
$ws.Range("K113").Value = 1679.8
$ws.Range("M113").Value = 490.2

# ---------------------------------------------------------------
# Sheet: LTW
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")
# Leve row 61
$ws.Range("H61").Value = 5639.4
$ws.Range("I61").Value = 2438.4
$ws.Range("K61").Value = 2438.4
$ws.Range("M61").Value = -2236.4
# Leve row 82
$ws.Range("H82").Value = 1589.375
$ws.Range("J82").Value = 1677.4445
$ws.Range("L82").Value = 1677.4445
$ws.Range("N82").Value = -2399.4445
# Leve row 85
$ws.Range("H85").Value = 1589.375
$ws.Range("J85").Value = 1677.4445
$ws.Range("L85").Value = 1677.4445
$ws.Range("N85").Value = -4173.4445
# Leve row 100
$ws.Range("H100").Value = 7115.909
$ws.Range("I100").Value = 2365.3845
$ws.Range("K100").Value = 2365.3845
$ws.Range("M100").Value = -1824.3845
# Leve row 113
$ws.Range("H113").Value = 5639.4
$ws.Range("I113").Value = 2438.4
$ws.Range("K113").Value = 2438.4
$ws.Range("M113").Value = -268.4000000000001
# Leve row 132
$ws.Range("H132").Value = 3712.2856
$ws.Range("I132").Value = 3409.7646
$ws.Range("K132").Value = 10229.2938
$ws.Range("M132").Value = -7699.293799999999
# Leve row 136
$ws.Range("H136").Value = 4210.1777
$ws.Range("I136").Value = 3786.805
$ws.Range("K136").Value = 11360.415
$ws.Range("M136").Value = -8810.414999999999

# ---------------------------------------------------------------
# Sheet: WVR
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")
# Leve row 81
$ws.Range("H81").Value = 981.13336
$ws.Range("I81").Value = 936
$ws.Range("J81").Value = 1011.2222
$ws.Range("K81").Value = 1872
$ws.Range("L81").Value = 2022.4444
$ws.Range("M81").Value = -811
$ws.Range("N81").Value = -4144.4444
# Leve row 84
$ws.Range("H84").Value = 981.13336
$ws.Range("I84").Value = 936
$ws.Range("J84").Value = 1011.2222
$ws.Range("K84").Value = 9360
$ws.Range("L84").Value = 10112.222
$ws.Range("M84").Value = -4056
$ws.Range("N84").Value = -20720.222
# Leve row 113
$ws.Range("H113").Value = 2611.96
$ws.Range("J113").Value = 4305.1816
$ws.Range("L113").Value = 12915.5448
$ws.Range("N113").Value = -17255.5448
# Leve row 122
$ws.Range("H122").Value = 2205.8635
$ws.Range("I122").Value = 1712.0526
$ws.Range("K122").Value = 5136.1578
$ws.Range("M122").Value = -2686.1578

